$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section "6. Korisniku blokiranog novcanika..." (rows 48-50)
# Fill in the "Realno utroseno vreme" (actual time spent) column for the
# three tasks of the "block/unblock wallet" user story.
$ws.Range("C48").Value = "10min"
$ws.Range("C49").Value = "15min"
$ws.Range("C50").Value = "10min"

# Section "9. Administrator moze da blokira/odblokira odredjeni novcanik/nalog" (rows 66-68)
$ws.Range("C66").Value = "10min"
$ws.Range("C67").Value = "20min"
$ws.Range("C68").Value = "10min"

# Reflect the scrolled/selected state shown in the saved workbook.
[void]$ws.Range("C68").Select()
